# Update workbook to 23 agosto 2021 - append new rows 344-357 to Sheet1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows to append: row, date(serial), B (nuovi pos.), C (somma mobile 7gg.), D (per 100mila abitanti)
$data = @(
    @(344, 44418, 1, 1, 43.78283712784589),
    @(345, 44419, 0, 1, 43.78283712784589),
    @(346, 44420, 0, 1, 43.78283712784589),
    @(347, 44421, 0, 1, 43.78283712784589),
    @(348, 44422, 2, 3, 131.3485113835376),
    @(349, 44423, 3, 6, 262.6970227670753),
    @(350, 44424, 0, 6, 262.6970227670753),
    @(351, 44425, 1, 6, 262.6970227670753),
    @(352, 44426, 0, 6, 262.6970227670753),
    @(353, 44427, 0, 6, 262.6970227670753),
    @(354, 44428, 0, 6, 262.6970227670753),
    @(355, 44429, 0, 4, 175.1313485113835),
    @(356, 44430, 1, 2, 87.56567425569177),
    @(357, 44431, 0, 2, 87.56567425569177)
)

foreach ($entry in $data) {
    $r = $entry[0]
    $dateSerial = $entry[1]
    $b = $entry[2]
    $c = $entry[3]
    $d = $entry[4]

    # copy cell formatting (style) from the row above, like the existing pattern
    $ws.Range("A" + ($r - 1)).Copy($ws.Range("A" + $r))

    $ws.Cells.Item($r, 1).Value = $dateSerial
    $ws.Cells.Item($r, 2).Value = $b
    $ws.Cells.Item($r, 3).Value = $c
    $ws.Cells.Item($r, 4).Value = $d
}
